# ---------------------------------------------------------------------------
# edit.ps1 - applies the three text edits described by the unified diff:
#
#   1. "...folded across chest" + " " + "(with crossed wrists)" -> the two
#      trailing runs (a space-only run and the parenthetical run) collapse
#      into a single run " (with crossed wrists)" (no visible text change,
#      pure run-merge).
#
#   2. "GO signal. Only the data coming from the force plates are used for
#      calculating this PI." -> "force plates" becomes "kinematics", with
#      the sentence split across three runs.
#
#   3. (STS CoP stability / PI3 paragraph) "Data from both the Chair and
#      lower limb kinematics are needed for calculating this PI." ->
#      "both the Chair and " is removed ("Data from lower limb kinematics
#      ..."), with the remainder split across four runs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- helpers ----------------------------------------------------------------

# Returns a Range for the Nth (1-based) occurrence of $searchText, searching
# the document body starting at character offset $searchStart.
function Get-NthMatchRange {
    param($doc, $searchText, $n, $searchStart)
    $pos = $searchStart
    $cnt = 0
    while ($true) {
        $r = $doc.Range($pos, $doc.Content.End)
        $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { return $null }
        $cnt = $cnt + 1
        if ($cnt -eq $n) {
            return $doc.Range($r.Start, $r.End)
        }
        $pos = $r.End
    }
}

# Forces a run boundary at [start, end) by toggling Bold on then off again,
# which makes the interop layer re-materialise that span as its own run(s)
# without leaving any visible/semantic formatting behind.
function Force-Split {
    param($doc, $start, $end)
    $r = $doc.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

# ===========================================================================
# Change 1: merge " " + "(with crossed wrists)" runs (first occurrence only -
# the one in the "5STS" protocol section; a later, already-merged duplicate
# of this same sentence exists in the "30sSTS" section and must stay as-is).
# ===========================================================================

$anchor1 = Get-NthMatchRange $d "cross chest (with crossed wrists). It is recommended" 1 0
if ($anchor1 -eq $null) { throw "change 1 anchor not found" }
$base1 = $anchor1.Start

$segCrossChest = "cross chest"
$segSpace = " "
$segParen = "(with crossed wrists)"
$segDot = "."
$segAfter = " It is recommended that the feet "

$p0 = $base1
$p1 = $p0 + $segCrossChest.Length
$p2 = $p1 + $segSpace.Length
$p3 = $p2 + $segParen.Length
$p4 = $p3 + $segDot.Length
$p5 = $p4 + $segAfter.Length

# sanity check the slice boundaries before mutating anything
$chk = $d.Range($p0, $p5)
if ($chk.Text -ne ($segCrossChest + $segSpace + $segParen + $segDot + $segAfter)) {
    throw "change 1 offsets do not line up: [$($chk.Text)]"
}

# Trigger an actual content edit over the " (with crossed wrists)" span so
# the engine recomputes/merges its runs, then restore the original text.
$mergeSpan = $d.Range($p1, $p3)
$originalMergeText = $mergeSpan.Text
$d.Range($p1, $p3).Text = "Z" + $originalMergeText.Substring(1)
$d.Range($p1, $p3).Text = $originalMergeText

# Re-assert the run boundaries that must remain untouched around the merge.
Force-Split $d $p0 $p1   # "cross chest"
Force-Split $d $p3 $p4   # "."
Force-Split $d $p4 $p5   # " It is recommended that the feet "

# ===========================================================================
# Change 2: "force plates" -> "kinematics" (PI1 / 30sSTS repetitions
# paragraph).
# ===========================================================================

$anchor2 = Get-NthMatchRange $d "executed in the 30s after the GO signal. Only the data coming from the force plates are used for calculating this PI." 1 0
if ($anchor2 -eq $null) { throw "change 2 anchor not found" }
$base2 = $anchor2.Start

$seg1 = "executed in the 30s after the "
$seg2 = "GO signal. Only the data coming from the "
$seg3 = "force plates"
$seg4 = " are used for calculating this PI."

$q0 = $base2
$q1 = $q0 + $seg1.Length
$q2 = $q1 + $seg2.Length
$q3 = $q2 + $seg3.Length
$q4 = $q3 + $seg4.Length

$chk2 = $d.Range($q0, $q4)
if ($chk2.Text -ne ($seg1 + $seg2 + $seg3 + $seg4)) {
    throw "change 2 offsets do not line up: [$($chk2.Text)]"
}

# Replace "force plates" with "kinematics".
$d.Range($q2, $q3).Text = "kinematics"
$newLen = "kinematics".Length
$q3new = $q2 + $newLen
$q4new = $q3new + $seg4.Length

# Force the three runs the diff expects, and keep the preceding
# "executed in the 30s after the " run from being swept into the merge.
Force-Split $d $q0 $q1          # "executed in the 30s after the "
Force-Split $d $q1 $q2          # "GO signal. Only the data coming from the "
Force-Split $d $q2 $q3new       # "kinematics"
Force-Split $d $q3new $q4new    # " are used for calculating this PI."

# ===========================================================================
# Change 3: STS CoP stability (PI3) paragraph - drop "both the Chair and ".
# ===========================================================================

$anchor3 = Get-NthMatchRange $d "STS CoP stability" 3 0
if ($anchor3 -eq $null) { throw "change 3 anchor not found" }

$fullAnchor = Get-NthMatchRange $d ". Data from both the Chair and lower limb kinematics are needed for calculating this PI. " 1 $anchor3.Start
if ($fullAnchor -eq $null) { throw "change 3 sentence not found" }
$base3 = $fullAnchor.Start

$t1 = ". "
$t2 = "D"
$t3 = "ata from "
$t4 = "both the Chair and "
$t5 = "lower limb kinematics are needed for calculating this PI. "

$z0 = $base3
$z1 = $z0 + $t1.Length
$z2 = $z1 + $t2.Length
$z3 = $z2 + $t3.Length
$z4 = $z3 + $t4.Length
$z5 = $z4 + $t5.Length

$chk3 = $d.Range($z0, $z5)
if ($chk3.Text -ne ($t1 + $t2 + $t3 + $t4 + $t5)) {
    throw "change 3 offsets do not line up: [$($chk3.Text)]"
}

# Delete "both the Chair and ".
$d.Range($z3, $z4).Text = ""
$z5new = $z3 + $t5.Length
$z3minus1 = $z3 - 1

# Force the four runs the diff expects: ". " | "D" | "ata from" | " lower
# limb kinematics are needed for calculating this PI. "
Force-Split $d $z0 $z1            # ". "
Force-Split $d $z1 $z2            # "D"
Force-Split $d $z2 $z3minus1      # "ata from"  (t3 minus its trailing space)
Force-Split $d $z3minus1 $z5new   # " lower limb kinematics are needed for calculating this PI. "

Write-Output "done"
